$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.722.40'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '3.403.32'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.13'
$ws.Range("E5").Value = '  -3.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.26'
$ws.Range("E6").Value = '  -3.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.36'
$ws.Range("E7").Value = '  -5.70%  '
$ws.Range("E8").Value = '  -4.17%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.954'
$ws.Range("E10").Value = '  -2.66%  '
$ws.Range("D11").Value = '3.403.04'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.55'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '4.047.20'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '92.596.34'
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.04'
$ws.Range("E18").Value = '  -4.25%  '
$ws.Range("D19").Value = '3.406.92'
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.78'
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.48'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.83'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.433'
$ws.Range("E24").Value = '  -12.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.51'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000183'
$ws.Range("E26").Value = '  -5.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.06'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.579.02'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.83'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.21'
$ws.Range("E31").Value = '  -4.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.71'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.134'
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("E35").Value = '  -4.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.67'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '548.88'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.40'
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.38'
$ws.Range("E41").Value = '  -5.15%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.149'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.911'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.64'
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.65'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.43'
$ws.Range("E47").Value = '  -4.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0403'
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.28'
$ws.Range("E49").Value = '  -3.53%  '
$ws.Range("E50").Value = '  -4.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("E51").Value = '  -1.69%  '
